# Refresh the crypto market snapshot (Price / Volume(1h) columns, and for the
# four re-ranked coin pairs also the Coin name + Link) to match the latest
# GitHub Actions scrape, per the commit's canonical diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '43.699.86'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.18%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.293.47'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.18%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '96.47'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.87%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '269.37'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("E7").Value = '  -0.80%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("E9").Value = '  -2.27%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '45.25'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.71%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0935'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.02%  '

# Row 12
$ws.Range("E12").Value = '  -3.18%  '

# Row 13
$ws.Range("E13").Value = '  +0.84%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '15.70'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.16%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.638.78'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -0.66%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.853'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.73%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.294.67'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.00%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '43.709.09'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.01%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.0000112'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +4.12%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.20'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.67%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '72.10'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.79%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '2.54'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +11.32%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '232.81'
$c.Style = "Normal"

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '9.12'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -5.33%  '

# Row 25
$ws.Range("E25").Value = '  +7.18%  '

# Row 26
$ws.Range("E26").Value = '  -0.11%  '

# Row 27
$ws.Range("E27").Value = '  -1.24%  '

# Row 28
$ws.Range("E28").Value = '  -1.67%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.28'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.49%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '38.47'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.61%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '174.76'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.72%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '21.84'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.70%  '

# Row 33
$ws.Range("E33").Value = '  +0.25%  '

# Row 34
$ws.Range("E34").Value = '  -2.24%  '

# Row 35
$ws.Range("E35").Value = '  -0.43%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.51'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.64%  '

# Row 37
$ws.Range("E37").Value = '  -1.96%  '

# Row 38
$ws.Range("E38").Value = '  -3.32%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.41'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.18%  '

# Row 40
$ws.Range("E40").Value = '  +1.34%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.10%  '

# Row 42
$ws.Range("E42").Value = '  -0.18%  '

# Row 43
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.34'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.74%  '

# Row 44
$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '64.50'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +4.03%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.76'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.02%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '5.17'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -5.59%  '

# Row 47
$ws.Range("E47").Value = '  -0.49%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '97.35'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.60%  '

# Row 49
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.20'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.98%  '

# Row 50
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.53'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +12.32%  '

# Row 51
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.434'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +3.12%  '
